$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.9007992734150422
$ws.Range("D2").Value = 1592676.245928357
$ws.Range("E2").Value = 1262.012775659722
$ws.Range("F2").Value = 820.8314352013583
$ws.Range("G2").Value = 0.434870620064634
$ws.Range("C3").Value = 0.9557531538107251
$ws.Range("D3").Value = 710386.9428068224
$ws.Range("E3").Value = 842.8445543555599
$ws.Range("F3").Value = 406.6482202447164
$ws.Range("G3").Value = 0.11311280279881
$ws.Range("C4").Value = 0.9602673466281921
$ws.Range("D4").Value = 637911.1866563606
$ws.Range("E4").Value = 798.6934246983386
$ws.Range("F4").Value = 438.7444573779662
$ws.Range("G4").Value = 0.1254996070529405
$ws.Range("C5").Value = 0.981202859910593
$ws.Range("D5").Value = 301789.7100395415
$ws.Range("E5").Value = 549.353902361257
$ws.Range("F5").Value = 274.0751989835267
$ws.Range("G5").Value = 0.06801154252413685
$ws.Range("C6").Value = 0.9808709916112536
$ws.Range("D6").Value = 307117.884291187
$ws.Range("E6").Value = 554.1821760857949
$ws.Range("F6").Value = 270.1861757215501
$ws.Range("G6").Value = 0.06541838138078189
$ws.Range("C7").Value = 0.9776175877457847
$ws.Range("D7").Value = 359351.564762329
$ws.Range("E7").Value = 599.4593937560149
$ws.Range("F7").Value = 318.0473490620614
$ws.Range("G7").Value = 0.09098810339732286
$ws.Range("C8").Value = 0.9723619202467083
$ws.Range("D8").Value = 443731.7610616764
$ws.Range("E8").Value = 666.1319396798779
$ws.Range("F8").Value = 433.1961704694404
$ws.Range("G8").Value = 0.1805299874145258
